$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number (45227 = 2023-10-28)
# that was bumped by one day (45228 = 2023-10-29) for every data row
# (rows 2 through 536).
$ws.Range("C2:C536").Value = 45228
